$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs, Col5a2, Ddr1, ECs)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.664368333333333
$ws.Range("H2").Value = 10.993105
$ws.Range("I2").Value = 0.01377800822483461
$ws.Range("J2").Value = 0.01377800822483461
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.2763116666666667
$ws.Range("N2").Value = 0.828935
$ws.Range("O2").Value = 0.02083107478128044
$ws.Range("P2").Value = 0.02083107478128044
$ws.Range("Q2").Value = 1.012507721463889
$ws.Range("R2").Value = 9.112569493175
$ws.Range("S2").Value = 0.0002870107196686268
$ws.Range("T2").Value = 0.0002870107196686268

# Row 3 (ECs, Col5a2, Ddr1, FAPs)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.664368333333333
$ws.Range("H3").Value = 10.993105
$ws.Range("I3").Value = 0.01377800822483461
$ws.Range("J3").Value = 0.01377800822483461
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.180798333333333
$ws.Range("N3").Value = 6.542395
$ws.Range("O3").Value = 0.1644098988384798
$ws.Range("P3").Value = 0.1644098988384798
$ws.Range("Q3").Value = 7.991248354052777
$ws.Range("R3").Value = 71.921235186475
$ws.Range("S3").Value = 0.002265240938440801
$ws.Range("T3").Value = 0.002265240938440801

# Row 4 (ECs, Col5a2, Ddr1, sCs)
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.664368333333333
$ws.Range("H4").Value = 10.993105
$ws.Range("I4").Value = 0.01377800822483461
$ws.Range("J4").Value = 0.01377800822483461
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 10.807288
$ws.Range("N4").Value = 32.421864
$ws.Range("O4").Value = 0.8147590263802398
$ws.Range("P4").Value = 0.8147590263802398
$ws.Range("Q4").Value = 39.60188391641334
$ws.Range("R4").Value = 356.41695524772
$ws.Range("S4").Value = 0.01122575656672519
$ws.Range("T4").Value = 0.01122575656672519

# Row 5 (FAPs, Col5a2, Ddr1, ECs)
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 218.9522706666667
$ws.Range("H5").Value = 656.856812
$ws.Range("I5").Value = 0.8232595393453117
$ws.Range("J5").Value = 0.8232595393453117
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.2763116666666667
$ws.Range("N5").Value = 0.828935
$ws.Range("O5").Value = 0.02083107478128044
$ws.Range("P5").Value = 0.02083107478128044
$ws.Range("Q5").Value = 60.49906682835778
$ws.Range("R5").Value = 544.49160145522
$ws.Range("S5").Value = 0.01714938102850467
$ws.Range("T5").Value = 0.01714938102850467

# Row 6 (FAPs, Col5a2, Ddr1, FAPs)
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 218.9522706666667
$ws.Range("H6").Value = 656.856812
$ws.Range("I6").Value = 0.8232595393453117
$ws.Range("J6").Value = 0.8232595393453117
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.180798333333333
$ws.Range("N6").Value = 6.542395
$ws.Range("O6").Value = 0.1644098988384798
$ws.Range("P6").Value = 0.1644098988384798
$ws.Range("Q6").Value = 477.4907469494155
$ws.Range("R6").Value = 4297.41672254474
$ws.Range("S6").Value = 0.1353520175815761
$ws.Range("T6").Value = 0.1353520175815761

# Row 7 (FAPs, Col5a2, Ddr1, sCs)
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 218.9522706666667
$ws.Range("H7").Value = 656.856812
$ws.Range("I7").Value = 0.8232595393453117
$ws.Range("J7").Value = 0.8232595393453117
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 10.807288
$ws.Range("N7").Value = 32.421864
$ws.Range("O7").Value = 0.8147590263802398
$ws.Range("P7").Value = 0.8147590263802398
$ws.Range("Q7").Value = 2366.280247348619
$ws.Range("R7").Value = 21296.52222613757
$ws.Range("S7").Value = 0.6707581407352309
$ws.Range("T7").Value = 0.6707581407352309

# Row 8 (sCs, Col5a2, Ddr1, ECs)
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 43.34113033333333
$ws.Range("H8").Value = 130.023391
$ws.Range("I8").Value = 0.1629624524298537
$ws.Range("J8").Value = 0.1629624524298537
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.2763116666666667
$ws.Range("N8").Value = 0.828935
$ws.Range("O8").Value = 0.02083107478128044
$ws.Range("P8").Value = 0.02083107478128044
$ws.Range("Q8").Value = 11.97565995762056
$ws.Range("R8").Value = 107.780939618585
$ws.Range("S8").Value = 0.003394683033107138
$ws.Range("T8").Value = 0.003394683033107139

# Row 9 (sCs, Col5a2, Ddr1, FAPs)
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 43.34113033333333
$ws.Range("H9").Value = 130.023391
$ws.Range("I9").Value = 0.1629624524298537
$ws.Range("J9").Value = 0.1629624524298537
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.180798333333333
$ws.Range("N9").Value = 6.542395
$ws.Range("O9").Value = 0.1644098988384798
$ws.Range("P9").Value = 0.1644098988384798
$ws.Range("Q9").Value = 94.51826479571611
$ws.Range("R9").Value = 850.664383161445
$ws.Range("S9").Value = 0.02679264031846281
$ws.Range("T9").Value = 0.02679264031846282

# Row 10 (sCs, Col5a2, Ddr1, sCs)
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 43.34113033333333
$ws.Range("H10").Value = 130.023391
$ws.Range("I10").Value = 0.1629624524298537
$ws.Range("J10").Value = 0.1629624524298537
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 10.807288
$ws.Range("N10").Value = 32.421864
$ws.Range("O10").Value = 0.8147590263802398
$ws.Range("P10").Value = 0.8147590263802398
$ws.Range("Q10").Value = 468.4000777578693
$ws.Range("R10").Value = 4215.600699820824
$ws.Range("S10").Value = 0.1327751290782837
$ws.Range("T10").Value = 0.1327751290782838
